$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.050234594880861
$ws.Range("D2").Value = 1.055046071215303
$ws.Range("E2").Value = 1.047268018390363
$ws.Range("F2").Value = 1.064593979995573
$ws.Range("I2").Value = 1.041238722882342
$ws.Range("J2").Value = 1.055268825570934
$ws.Range("K2").Value = 1.057787607566292
$ws.Range("L2").Value = 1.050031131613104
$ws.Range("M2").Value = 1.067309496387584
$ws.Range("N2").Value = 1.056767428344271

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.051848024733527
$ws.Range("D3").Value = 1.056302580382626
$ws.Range("E3").Value = 1.048664479952349
$ws.Range("F3").Value = 1.065970914027681
$ws.Range("I3").Value = 1.041622705300265
$ws.Range("J3").Value = 1.056528613046005
$ws.Range("K3").Value = 1.058856349436392
$ws.Range("L3").Value = 1.051237887165316
$ws.Range("M3").Value = 1.068500259838711
$ws.Range("N3").Value = 1.058029004862056

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.052890266349497
$ws.Range("D4").Value = 1.057113884361921
$ws.Range("E4").Value = 1.049566676595772
$ws.Range("F4").Value = 1.066860420077344
$ws.Range("I4").Value = 1.041869011810759
$ws.Range("J4").Value = 1.057341676636966
$ws.Range("K4").Value = 1.05954560124387
$ws.Range("L4").Value = 1.052016809856728
$ws.Range("M4").Value = 1.069268767928383
$ws.Range("N4").Value = 1.058843223096576

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.053328014200277
$ws.Range("D5").Value = 1.057454546151219
$ws.Range("E5").Value = 1.049945630588313
$ws.Range("F5").Value = 1.067234025201292
$ws.Range("I5").Value = 1.041972045486628
$ws.Range("J5").Value = 1.057682992379208
$ws.Range("K5").Value = 1.059834819196454
$ws.Range("L5").Value = 1.052343814036757
$ws.Range("M5").Value = 1.069591377348242
$ws.Range("N5").Value = 1.059185023546321

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.053401490198511
$ws.Range("D6").Value = 1.057511720886615
$ws.Range("E6").Value = 1.050009239526047
$ws.Range("F6").Value = 1.067296735217689
$ws.Range("I6").Value = 1.041989315250362
$ws.Range("J6").Value = 1.057740271908653
$ws.Range("K6").Value = 1.059883348467499
$ws.Range("L6").Value = 1.052398692998949
$ws.Range("M6").Value = 1.069645517430833
$ws.Range("N6").Value = 1.059242384419268

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.052896117157264
$ws.Range("D7").Value = 1.057118437904874
$ws.Range("E7").Value = 1.049571741481751
$ws.Range("F7").Value = 1.066865413546249
$ws.Range("I7").Value = 1.041870390566231
$ws.Range("J7").Value = 1.057346239254092
$ws.Range("K7").Value = 1.059549467915387
$ws.Range("L7").Value = 1.052021181080517
$ws.Range("M7").Value = 1.069273080494526
$ws.Range("N7").Value = 1.058847792193142

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.050780230384715
$ws.Range("D8").Value = 1.055471078885702
$ws.Range("E8").Value = 1.047740255016276
$ws.Range("F8").Value = 1.065059628234114
$ws.Range("I8").Value = 1.04136893943115
$ws.Range("J8").Value = 1.055695016528171
$ws.Range("K8").Value = 1.058149272883501
$ws.Range("L8").Value = 1.050439363818281
$ws.Range("M8").Value = 1.067712337480711
$ws.Range("N8").Value = 1.057194224541552

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.047037886446619
$ws.Range("D9").Value = 1.052554591122467
$ws.Range("E9").Value = 1.04450182790506
$ws.Range("F9").Value = 1.061866094761165
$ws.Range("I9").Value = 1.04046869406247
$ws.Range("J9").Value = 1.052768914763738
$ws.Range("K9").Value = 1.055664094829316
$ws.Range("L9").Value = 1.047636918789448
$ws.Range("M9").Value = 1.06494654019972
$ws.Range("N9").Value = 1.054263967376985

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.044533044456375
$ws.Range("D10").Value = 1.050600688886368
$ws.Range("E10").Value = 1.042334932936206
$ws.Range("F10").Value = 1.059728899226115
$ws.Range("I10").Value = 1.039857192816268
$ws.Range("J10").Value = 1.050806667263452
$ws.Range("K10").Value = 1.053994919374246
$ws.Range("L10").Value = 1.045758049190206
$ws.Range("M10").Value = 1.063091805798621
$ws.Range("N10").Value = 1.052298933260176

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.043445924086965
$ws.Range("D11").Value = 1.049752263488124
$ws.Range("E11").Value = 1.041394653826493
$ws.Range("F11").Value = 1.058801433721245
$ws.Range("I11").Value = 1.039589680741166
$ws.Range("J11").Value = 1.049954158166301
$ws.Range("K11").Value = 1.053269122788169
$ws.Range("L11").Value = 1.044941873040758
$ws.Range("M11").Value = 1.06228601696209
$ws.Range("N11").Value = 1.051445213502322

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.043041730605073
$ws.Range("D12").Value = 1.049436755953828
$ws.Range("E12").Value = 1.041045082751518
$ws.Range("F12").Value = 1.058456615880131
$ws.Range("I12").Value = 1.039489901813548
$ws.Range("J12").Value = 1.049637062194951
$ws.Range("K12").Value = 1.052999066466246
$ws.Range("L12").Value = 1.044638307820684
$ws.Range("M12").Value = 1.061986301109494
$ws.Range("N12").Value = 1.051127667218322

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.043128449254648
$ws.Range("D13").Value = 1.049504450018713
$ws.Range("E13").Value = 1.041120081089643
$ws.Range("F13").Value = 1.058530594908205
$ws.Range("I13").Value = 1.039511323474186
$ws.Range("J13").Value = 1.049705100238722
$ws.Range("K13").Value = 1.053057015544502
$ws.Range("L13").Value = 1.044703441876083
$ws.Range("M13").Value = 1.062050609857847
$ws.Range("N13").Value = 1.051195801883919

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.043412521323545
$ws.Range("D14").Value = 1.049726191020906
$ws.Range("E14").Value = 1.041365764540443
$ws.Range("F14").Value = 1.05877293746382
$ws.Range("I14").Value = 1.039581441428117
$ws.Range("J14").Value = 1.049927955863807
$ws.Range("K14").Value = 1.053246809340415
$ws.Range("L14").Value = 1.044916788475591
$ws.Range("M14").Value = 1.062261250765908
$ws.Range("N14").Value = 1.051418973989554

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.043587495653797
$ws.Range("D15").Value = 1.049862764419499
$ws.Range("E15").Value = 1.041517096929048
$ws.Range("F15").Value = 1.058922210600631
$ws.Range("I15").Value = 1.039624588587302
$ws.Range("J15").Value = 1.05006520650662
$ws.Range("K15").Value = 1.053363685962941
$ws.Range("L15").Value = 1.045048184958668
$ws.Range("M15").Value = 1.062390979015553
$ws.Range("N15").Value = 1.051556419544022

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.044605139217731
$ws.Range("D16").Value = 1.050656945405485
$ws.Range("E16").Value = 1.042397293259799
$ws.Range("F16").Value = 1.059790408231366
$ws.Range("I16").Value = 1.039874888965031
$ws.Range("J16").Value = 1.050863184819083
$ws.Range("K16").Value = 1.054043023567927
$ws.Range("L16").Value = 1.045812160336099
$ws.Range("M16").Value = 1.063145226286724
$ws.Range("N16").Value = 1.052355531077218

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.045242800516129
$ws.Range("D17").Value = 1.051154473203779
$ws.Range("E17").Value = 1.042948875443516
$ws.Range("F17").Value = 1.060334451983399
$ws.Range("I17").Value = 1.040031163179848
$ws.Range("J17").Value = 1.051362968050216
$ws.Range("K17").Value = 1.054468336885106
$ws.Range("L17").Value = 1.046290676331079
$ws.Range("M17").Value = 1.063617623535763
$ws.Range("N17").Value = 1.052856024057865

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.045614496457391
$ws.Range("D18").Value = 1.051444444356464
$ws.Range("E18").Value = 1.0432704117214
$ws.Range("F18").Value = 1.060651586695044
$ws.Range("I18").Value = 1.040122052289655
$ws.Range("J18").Value = 1.051654209459934
$ws.Range("K18").Value = 1.054716122858782
$ws.Range("L18").Value = 1.046569535053008
$ws.Range("M18").Value = 1.063892907141546
$ws.Range("N18").Value = 1.053147679063791

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.045741194630028
$ws.Range("D19").Value = 1.051543278529197
$ws.Range("E19").Value = 1.043380014945445
$ws.Range("F19").Value = 1.060759688310752
$ws.Range("I19").Value = 1.040152998614314
$ws.Range("J19").Value = 1.051753469102292
$ws.Range("K19").Value = 1.054800562156014
$ws.Range("L19").Value = 1.046664576255898
$ws.Range("M19").Value = 1.063986728260651
$ws.Range("N19").Value = 1.053247079666227

$ws.Range("B20").Value = 1.019999999999999
$ws.Range("C20").Value = 1.045174410518281
$ws.Range("D20").Value = 1.051101116861959
$ws.Range("E20").Value = 1.042889715849966
$ws.Range("F20").Value = 1.060276101655609
$ws.Range("I20").Value = 1.040014423661484
$ws.Range("J20").Value = 1.051309374414015
$ws.Range("K20").Value = 1.054422735057971
$ws.Range("L20").Value = 1.04623936217488
$ws.Range("M20").Value = 1.06356696649557
$ws.Range("N20").Value = 1.052802354312554

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.043328880002032
$ws.Range("D21").Value = 1.049660903951503
$ws.Range("E21").Value = 1.04129342550682
$ws.Range("F21").Value = 1.058701582409346
$ws.Range("I21").Value = 1.039560804865728
$ws.Range("J21").Value = 1.049862342548953
$ws.Range("K21").Value = 1.053190932644213
$ws.Range("L21").Value = 1.044853974358974
$ws.Range("M21").Value = 1.062199233650183
$ws.Range("N21").Value = 1.051353267496266

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.042166267594911
$ws.Range("D22").Value = 1.048753271617536
$ws.Range("E22").Value = 1.040287978859205
$ws.Range("F22").Value = 1.057709787943378
$ws.Range("I22").Value = 1.039273205451707
$ws.Range("J22").Value = 1.048950008321402
$ws.Range("K22").Value = 1.052413766126721
$ws.Range("L22").Value = 1.043980601414639
$ws.Range("M22").Value = 1.061336909634376
$ws.Range("N22").Value = 1.050439637649466

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.042782808087952
$ws.Range("D23").Value = 1.049234627722804
$ws.Range("E23").Value = 1.04082115843281
$ws.Range("F23").Value = 1.058235733441618
$ws.Range("I23").Value = 1.039425895077247
$ws.Range("J23").Value = 1.049433896597409
$ws.Range("K23").Value = 1.052826013696927
$ws.Range("L23").Value = 1.044443816126847
$ws.Range("M23").Value = 1.061794271783112
$ws.Range("N23").Value = 1.050924213102328

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.045205313759596
$ws.Range("D24").Value = 1.051125227002039
$ws.Range("E24").Value = 1.042916448118089
$ws.Range("F24").Value = 1.060302468265398
$ws.Range("I24").Value = 1.040021988342003
$ws.Range("J24").Value = 1.051333591917129
$ws.Range("K24").Value = 1.054443341463421
$ws.Range("L24").Value = 1.04626254961072
$ws.Range("M24").Value = 1.063589857028029
$ws.Range("N24").Value = 1.0528266062073

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.048007079868412
$ws.Range("D25").Value = 1.05331023067313
$ws.Range("E25").Value = 1.045340404265102
$ws.Range("F25").Value = 1.062693108164003
$ws.Range("I25").Value = 1.040703415561613
$ws.Range("J25").Value = 1.053527377337078
$ws.Range("K25").Value = 1.056308728328137
$ws.Range("L25").Value = 1.048363249560071
$ws.Range("M25").Value = 1.069591377348242
$ws.Range("N25").Value = 1.055023507054176
